# StatusTableOverview - confirm code edits for the Train-Bridge Controller row.
#
# Row 3 (A3 = "Train-Bridge Controller") is updated to reflect that the code
# edits were confirmed:
#   - "Code Edits" (E3) comment is updated from "need to be confirmed" to
#     "removed comments", and its status highlight changes from the
#     "Neutral" (yellow) style to the "Good" (green) style.
#   - "Reached States" (F3) and "Reached Transitions" (G3) also switch from
#     "Neutral" to "Good" to reflect the confirmed status.
#
# Row 2's F2/G2 ("Reached States"/"Reached Transitions") are likewise
# switched from "Neutral" to "Good".
#
# Finally, the active selection moves to G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatusTableOverview")

# Row 2 ("Philosophers"): mark reached states/transitions as confirmed (Good).
$ws.Range("F2").Style = "Good"
$ws.Range("G2").Style = "Good"

# Row 3 ("Train-Bridge Controller"): update the code-edit comment and mark
# the row's status cells as confirmed (Good).
$ws.Range("E3").Value = "removed comments"
$ws.Range("E3").Style = "Good"
$ws.Range("F3").Style = "Good"
$ws.Range("G3").Style = "Good"

# Update the sheet's active selection to G3.
$ws.Range("G3").Select()
